$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 8 first (grandes regiões header row), then row 5 (situação do domicílio header row)
# Deleting from the bottom up keeps earlier row indices stable.
$ws.Rows("8:8").Delete()
$ws.Rows("5:5").Delete()
